# #5: property boat&car done
# Extend the "汽車" (car/boat) sheet (sheet index 3) with a new "capacity"
# column in the header row, and populate the standard trailer columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) on the data row, matching the layout already used
# by every other sheet in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Header row (row 1) -----------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row (row 2) ---------------------------------------------------
$ws.Range("A2").Value = 44
$ws.Range("B2").Value = "HondaCRV"
$ws.Range("C2").Value = 1997
$ws.Range("D2").Value = "黃玉廷"
$ws.Range("E2").Value = "96年07月11曰"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = 800000
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-24"
$ws.Range("K2").Value = "蔡其昌"
$ws.Range("L2").Value = 1377
$ws.Range("M2").Value = "tmp30a51"
$ws.Range("N2").Value = 44
